$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new values (were row 4's data, with Q/R rounded) ---
$ws.Range("A3").Value = 111714275
$ws.Range("B3").Value = 90660
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 4362
$ws.Range("F3").Value = "Blå taggsvamp"
$ws.Range("G3").Value = "Hydnellum caeruleum"
$ws.Range("H3").Value = "(Hornem.) P.Karst."
$ws.Range("Q3").Value = 513092
$ws.Range("R3").Value = 6552092
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Row 4: new values (were row 3's data, with Q/R rounded) ---
$ws.Range("A4").Value = 111714277
$ws.Range("B4").Value = 90671
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 6003298
$ws.Range("F4").Value = "Ruttaggsvamp"
$ws.Range("G4").Value = "Hydnellum illudens"
$ws.Range("H4").Value = "(Maas Geest.) Nitare"
$ws.Range("Q4").Value = 513092
$ws.Range("R4").Value = 6552092
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# --- Row 5: Q/R rounded, Starttid/Sluttid cleared ---
$ws.Range("Q5").Value = 513071
$ws.Range("R5").Value = 6552108
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
